$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column cells are stored as exact text (avoid numeric coercion
# that would round/reformat values like "163.00" -> 163 or "0.0678" -> 6.78E-02).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.554.80"
$ws.Range("E2").Value = "  +3.89%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.254.44"
$ws.Range("E3").Value = "  +3.84%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.38"
$ws.Range("E5").Value = "  +1.97%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "181.38"
$ws.Range("E6").Value = "  +7.18%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("E8").Value = "  +1.09%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.253.03"
$ws.Range("E9").Value = "  +3.98%  "

$ws.Range("E10").Value = "  +7.97%  "

$ws.Range("E11").Value = "  +3.11%  "

$ws.Range("E12").Value = "  +7.19%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.821.13"
$ws.Range("E13").Value = "  +4.13%  "

$ws.Range("E14").Value = "  +1.49%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "28.35"
$ws.Range("E15").Value = "  +5.29%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.518.59"
$ws.Range("E16").Value = "  +3.88%  "

$ws.Range("E17").Value = "  +4.12%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.256.38"
$ws.Range("E18").Value = "  +4.15%  "

$ws.Range("E19").Value = "  +3.37%  "

$ws.Range("E20").Value = "  +6.39%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "374.74"
$ws.Range("E21").Value = "  +5.32%  "

$ws.Range("E22").Value = "  +5.89%  "

$ws.Range("E23").Value = "  +0.08%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.06"
$ws.Range("E24").Value = "  +3.74%  "

$ws.Range("E25").Value = "  +4.28%  "

$ws.Range("E26").Value = "  +6.83%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.56"
$ws.Range("E27").Value = "  -0.40%  "

$ws.Range("E28").Value = "  +3.78%  "

$ws.Range("E29").Value = "  +0.14%  "

$ws.Range("E30").Value = "  +3.99%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.66"
$ws.Range("E31").Value = "  +8.11%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "22.70"
$ws.Range("E32").Value = "  +4.79%  "

$ws.Range("E33").Value = "  -0.02%  "

$ws.Range("E34").Value = "  +7.75%  "

$ws.Range("E35").Value = "  +6.12%  "

$ws.Range("E36").Value = "  +6.21%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "163.00"
$ws.Range("E37").Value = "  +2.68%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.847"
$ws.Range("E38").Value = "  +3.41%  "

$ws.Range("E39").Value = "  +5.88%  "

$ws.Range("E40").Value = "  +13.01%  "

$ws.Range("E41").Value = "  +2.80%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.63"
$ws.Range("E42").Value = "  +12.19%  "

$ws.Range("E43").Value = "  +7.19%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.696.47"
$ws.Range("E44").Value = "  +2.54%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "351.88"
$ws.Range("E45").Value = "  +9.74%  "

$ws.Range("E46").Value = "  +7.18%  "

$ws.Range("E47").Value = "  +3.23%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0678"
$ws.Range("E48").Value = "  +4.44%  "

$ws.Range("E49").Value = "  +3.64%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.995"
$ws.Range("E50").Value = "  +6.80%  "

$ws.Range("E51").Value = "  +1.05%  "
